# Weekly canteen-menu refresh published by the controller: new serving date
# (07.04.2025 -> 09.04.2025), new dish names and new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Serving date (column A) -------------------------------------------------
# A2:A10 all share the same literal date text. Writing the new date straight
# into a Range.Value makes the host re-interpret "09.04.2025" as an actual
# date serial (since those cells already carry a date display format), which
# is not what the source file stores (it stores the date as plain text).
# Build the literal text in a scratch cell via a formula -- a formula result
# is always treated as text/number verbatim, never re-parsed -- then paste
# just the *value* into each date cell so every cell keeps its own original
# formatting untouched.
$ws.Range("F1").Formula = '=""&"09.04.2025"'
$ws.Range("F1").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("A3").PasteSpecial(-4163)
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("A7").PasteSpecial(-4163)
$ws.Range("A8").PasteSpecial(-4163)
$ws.Range("A9").PasteSpecial(-4163)
$ws.Range("A10").PasteSpecial(-4163)
$ws.Range("F1").Clear()
$excel.CutCopyMode = $false

# A3 additionally picks up the same "date" cell styling its neighbours
# (A2, A5, A7:A9) already use -- copy the format only, leaving the value
# (already correct from the paste above) untouched.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Dish names (column B) ---------------------------------------------------
$ws.Range("B2").Value = "Schabik pod pierzynk" + [char]0x0105 + " majonezowo-warzywn" + [char]0x0105 + ", ziemniaki, sur" + [char]0x00F3 + "wka + zupa krupnik lub go" + [char]0x0142 + [char]0x0105 + "bkowa"
$kotletyPozarskie = "Kotlety po" + [char]0x017C + "arskie, ziemniaki, sur" + [char]0x00F3 + "wka + zupa krupnik lub go" + [char]0x0142 + [char]0x0105 + "bkowa"
$ws.Range("B3").Value = $kotletyPozarskie
$ws.Range("B4").Value = "Schabik pod pierzynk" + [char]0x0105 + ", ziemniaki, sur" + [char]0x00F3 + "wka"
$ws.Range("B5").Value = $kotletyPozarskie
$ws.Range("B6").Value = "Zupa krupnik lub go" + [char]0x0142 + [char]0x0105 + "bkowa"
$ws.Range("B9").Value = "Zupa gulaszowa"

# --- Prices (column C) --------------------------------------------------------
$ws.Range("C2").Value = 33
$ws.Range("C3").Value = 30
$ws.Range("C4").Value = 30
$ws.Range("C5").Value = 27
$ws.Range("C9").Value = 16.5
$ws.Range("C10").Value = 22
